$wb = $excel.ActiveWorkbook

# --- Template sheet ---
$ws1 = $wb.Worksheets.Item("Template")
$ws1.Columns("L").Insert()
$ws1.Range("L2").Value = "Group"
$ws1.Range("L2").Select()

# --- Sample Data sheet ---
$ws2 = $wb.Worksheets.Item("Sample Data")
$ws2.Columns("L").Insert()
$ws2.Range("L2").Value = "Group"
$ws2.Range("L3").Value = "Odell"
$ws2.Range("L4").Select()
